# "Generate Report for handback" — stamp the handoff/handback timestamps
# for the 9c62d602-...-776282620988 file on both locale report sheets
# (zh-cn and de-de), row 3 of each table.
#
#   Column D = Correspond Handoff Datetime
#   Column G = Correspond Handback DateTime

$wb = $excel.ActiveWorkbook

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-14 15:25:00"
$wsZhCn.Range("G3").Value = "2016-01-14 15:26:49"

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-14 15:25:45"
$wsDeDe.Range("G3").Value = "2016-01-14 15:27:43"
